# Team_Members.pptx — "Fixed Yong Wei's UID"
#
# 1) The slide master and every slide layout cache a fixed
#    "datetimeFigureOut" footer date ("7/08/2018"). Update the cached
#    text on the Date placeholder (shape index 3) on the master and on
#    every custom layout to "21/08/2018".
# 2) Yong Wei Lim's e-mail address text box on slide 1 shows the wrong
#    university UID. Update the run text from "u6019071@anu.edu.au" to
#    "u6033740@anu.edu.au" (slide 2's copy of this textbox already uses
#    a different hyperlink relationship id and is left untouched).

$p = $ppt.ActivePresentation

# --- 1) Fix the cached footer date on the master + all layouts -------
$master = $p.Slides.Item(1).Master

$masterDate = $master.Shapes.Item(3)
if ($masterDate.Name -eq "Date Placeholder 3") {
    $masterDate.TextFrame.TextRange.Text = "21/08/2018"
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    $layoutDate = $layout.Shapes.Item(3)
    if ($layoutDate.Name -eq "Date Placeholder 3") {
        $layoutDate.TextFrame.TextRange.Text = "21/08/2018"
    }
}

# --- 2) Fix Yong Wei's UID on slide 1 ---------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "u6019071@anu.edu.au") {
                $shp.TextFrame.TextRange.Text = "u6033740@anu.edu.au"
            }
        }
    }
}
